# BubbleSheet/answers.xlsx - "separate code in classes"
#
# The scanner app appended a fresh batch of 11 graded answer sheets
# (rows 67-77) to the "bubble sheet" tab, and the placeholder blank
# name cell it used to stamp on the previously-last row (B66) is gone
# now that the row no longer sits at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bubble sheet")

# B66 was an empty recognized-name placeholder; it is not re-emitted now
# that row 66 is no longer the last row in the sheet.
$ws.Range("B66").ClearContents()

# Column layout per scanned row: A = sheet/student id (text, may have
# leading zeros), B = recognized name (text, optional), C.. = one 0/1
# flag per bubble option.
$newRows = @(
    @{ Row = 67; Id = "2303";  Name = $null;           Answers = @(0,1,0,0,0,0,0,0,0,0,0,0,1,0) },
    @{ Row = 68; Id = "10301"; Name = "MohamedAhied";   Answers = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1) },
    @{ Row = 69; Id = "2303";  Name = $null;           Answers = @(0,1,0,0,0,0,0,0,0,0,0,0,1,0) },
    @{ Row = 70; Id = "2303";  Name = $null;           Answers = @(0,1,0,0,0,0,0,0,0,0,0,0,1,0) },
    @{ Row = 71; Id = "10301"; Name = "MohamedAhied";   Answers = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1) },
    @{ Row = 72; Id = "100";   Name = $null;           Answers = @(0,0,0,0,0,0,0,0,0,0,1,0,0,0) },
    @{ Row = 73; Id = "09211"; Name = $null;           Answers = @(0,1,0,0,1,0,0,0,0,0,0,1,0,0) },
    @{ Row = 74; Id = "2303";  Name = $null;           Answers = @(0,1,0,0,0,0,0,0,0,0,0,0,1,0) },
    @{ Row = 75; Id = "10301"; Name = "MohamedAhied";   Answers = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1) },
    @{ Row = 76; Id = "100";   Name = $null;           Answers = @(0,0,0,0,0,0,0,0,0,0,1,0,0,0) },
    @{ Row = 77; Id = "09211"; Name = "";               Answers = @(0,1,0,0,1,0,0,0,0,0,0,1,0,0) }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: force text so ids like "09211"/"000" keep their leading
    # zeros instead of being parsed as numbers.
    $ws.Cells.Item($row, 1).Value = "'" + $r.Id

    if ($null -ne $r.Name) {
        if ($r.Name -eq "") {
            # Row 77 mirrors row 66's old pattern: a recognized-name
            # column present on the row, just empty.
            $ws.Cells.Item($row, 2).NumberFormat = "@"
            $ws.Cells.Item($row, 2).Value = ""
        } else {
            $ws.Cells.Item($row, 2).Value = $r.Name
        }
    }

    # Columns C..Q: the per-option bubble flags (0/1).
    $col = 3
    foreach ($ans in $r.Answers) {
        $ws.Cells.Item($row, $col).Value = $ans
        $col++
    }
}
